$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 7 (old second-to-last games), keeping only header row 1 and data row 2
$ws.Range("A3:J7").EntireRow.Delete()

# Remove the "Home Bet" / "Away Bet" columns (I and J)
$ws.Range("I1:J1").EntireColumn.Delete()

# Update row 2 with the new game data
$ws.Range("A2").Value = "TOR"
$ws.Range("B2").Value = "NYK"
$ws.Range("C2").Value = 200
$ws.Range("D2").Value = -245
$ws.Range("E2").Value = 0.3507480227274955
$ws.Range("F2").Value = 0.6495308864939859
$ws.Range("G2").Value = -9.775593181751347
$ws.Range("H2").Value = -45.45454929892449
